# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.268.43'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '2.047.31'
$ws.Range("E3").Value = '  -1.57%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''229.68'
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").Value = '''0.618'
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''56.89'
$ws.Range("E8").Value = '  -3.26%  '
$ws.Range("D9").Value = '''0.382'
$ws.Range("E9").Value = '  -2.69%  '
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").Value = '''14.65'
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.348.80'
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").Value = '''20.42'
$ws.Range("E14").Value = '  -3.09%  '
$ws.Range("D15").Value = '''0.751'
$ws.Range("E15").Value = '  -3.15%  '
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").Value = '2.060.50'
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").Value = '37.221.02'
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").Value = '''69.47'
$ws.Range("E20").Value = '  -2.88%  '
$ws.Range("D21").Value = '0.0₃0819'
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").Value = '''225.44'
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '''2.34'
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  -3.97%  '
$ws.Range("D26").Value = '''9.46'
$ws.Range("E26").Value = '  +2.56%  '
$ws.Range("D27").Value = '''168.16'
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("E28").Value = '  -5.29%  '
$ws.Range("D29").Value = '''19.07'
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("D30").Value = '''1.34'
$ws.Range("E30").Value = '  -5.40%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").Value = '''4.51'
$ws.Range("E32").Value = '  -3.97%  '
$ws.Range("D33").Value = '''0.0621'
$ws.Range("E33").Value = '  -2.17%  '
$ws.Range("D34").Value = '''4.57'
$ws.Range("E34").Value = '  -2.74%  '
$ws.Range("D35").Value = '''2.46'
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("D37").Value = '''3.27'
$ws.Range("E37").Value = '  -3.79%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").Value = '  -2.61%  '
$ws.Range("D40").Value = '''0.0225'
$ws.Range("E40").Value = '  +4.35%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.488.87'
$ws.Range("E41").Value = '  +3.43%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''97.55'
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("D44").Value = '''0.0946'
$ws.Range("E44").Value = '  -3.29%  '
$ws.Range("E45").Value = '  +2.63%  '
$ws.Range("D46").Value = '''16.54'
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").Value = '''4.04'
$ws.Range("E47").Value = '  -3.76%  '
$ws.Range("E48").Value = '  -3.54%  '
$ws.Range("D49").Value = '''7.16'
$ws.Range("E49").Value = '  -3.49%  '
$ws.Range("E50").Value = '  -1.38%  '
$ws.Range("D51").Value = '2.234.48'
$ws.Range("E51").Value = '  -1.58%  '
